$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new date-serial (column A), new percent value (column B)
$rows = @(
    ,(2, 45785.375, 0.002000890919897798)
    ,(3, 45785.37847222222, 0.005146220767290485)
    ,(4, 45785.38194444445, 0.007309114207060435)
    ,(5, 45785.38541666666, 0.0528345396846035)
    ,(6, 45785.38888888889, 0.07682729025092071)
    ,(7, 45785.39236111111, 0.105587160459568)
    ,(8, 45785.39583333334, 0.1320499672136994)
    ,(9, 45785.39930555555, 0.1550265514964979)
    ,(10, 45785.40277777778, 0.1769693118514202)
    ,(11, 45785.40625, 0.1989684840795836)
    ,(12, 45785.40972222222, 0.2171710197813316)
    ,(13, 45785.41319444445, 0.232682418585247)
    ,(14, 45785.41666666666, 0.2445766829652835)
    ,(15, 45785.42013888889, 0.2602061288480008)
    ,(16, 45785.42361111111, 0.2733402915799901)
    ,(17, 45785.42708333334, 0.2899181318025941)
    ,(18, 45785.43055555555, 0.3057263198382481)
    ,(19, 45785.43402777778, 0.3217240613408628)
    ,(20, 45785.4375, 0.3367974953216982)
    ,(21, 45785.44097222222, 0.3508931795126293)
    ,(22, 45785.44444444445, 0.3655014648706726)
    ,(23, 45785.44791666666, 0.3798428831291147)
    ,(24, 45785.45138888889, 0.3931052458931673)
    ,(25, 45785.45486111111, 0.4072928831411613)
    ,(26, 45785.45833333334, 0.4193618208296876)
    ,(27, 45785.46180555555, 0.4319521053632188)
    ,(28, 45785.46527777778, 0.4443217043052283)
    ,(29, 45785.46875, 0.4559775489999275)
    ,(30, 45785.47222222222, 0.4681817390733403)
    ,(31, 45785.47569444445, 0.4817875831033827)
    ,(32, 45785.47916666666, 0.4939852457220174)
    ,(33, 45785.54166666666, 0.5082356443978371)
    ,(34, 45785.54513888889, 0.5293149471409798)
    ,(35, 45785.54861111111, 0.5491075733607118)
    ,(36, 45785.55208333334, 0.5706934972985998)
    ,(37, 45785.55555555555, 0.5931401138147578)
    ,(38, 45785.55902777778, 0.6157300390772882)
    ,(39, 45785.5625, 0.6378314097499402)
    ,(40, 45785.56597222222, 0.6602726602467011)
    ,(41, 45785.56944444445, 0.6827974516914428)
    ,(42, 45785.57291666666, 0.7059592475232496)
    ,(43, 45785.57638888889, 0.7265784670215074)
    ,(44, 45785.57986111111, 0.7490170380406582)
    ,(45, 45785.58333333334, 0.7713136563251427)
    ,(46, 45785.58680555555, 0.7961202444168476)
    ,(47, 45785.59027777778, 0.8204249891797251)
    ,(48, 45785.59375, 0.8455421383400709)
    ,(49, 45785.59722222222, 0.8733570072484979)
    ,(50, 45785.60069444445, 0.9037553045706305)
    ,(51, 45785.60416666666, 0.9347689877083782)
    ,(52, 45785.60763888889, 0.9377147464548174)
    ,(53, 45785.61111111111, 0.9400047401639366)
    ,(54, 45785.61458333334, 0.9882157395613993)
    ,(55, 45785.61805555555, 0.9931331633185496)
    ,(56, 45785.62152777778, 0.9973332952849573)
    ,(57, 45785.625, 1.0)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $dateVal = $entry[1]
    $pctVal = $entry[2]
    $ws.Cells.Item($r, 1).Value = $dateVal
    $ws.Cells.Item($r, 2).Value = $pctVal
}